# Auto-generated script to apply numeric updates to Moogle_Profits workbook sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 1161.7333
$ws.Cells.Item(4, 9).Value = 760.5833
$ws.Cells.Item(4, 10).Value = 2766.3333
$ws.Cells.Item(4, 11).Value = 760.5833
$ws.Cells.Item(4, 12).Value = 2766.3333
$ws.Cells.Item(4, 13).Value = -646.5833
$ws.Cells.Item(4, 14).Value = -2994.3333
$ws.Cells.Item(18, 8).Value = 772.8333
$ws.Cells.Item(18, 9).Value = 772.8333
$ws.Cells.Item(18, 11).Value = 772.8333
$ws.Cells.Item(18, 13).Value = -488.8333
$ws.Cells.Item(53, 8).Value = 275
$ws.Cells.Item(53, 10).Value = 301.125
$ws.Cells.Item(53, 12).Value = 301.125
$ws.Cells.Item(53, 14).Value = -1575.125
$ws.Cells.Item(62, 8).Value = 5601.7407
$ws.Cells.Item(62, 9).Value = 4135.65
$ws.Cells.Item(62, 11).Value = 4135.65
$ws.Cells.Item(62, 13).Value = -3511.65
$ws.Cells.Item(65, 8).Value = 5601.7407
$ws.Cells.Item(65, 9).Value = 4135.65
$ws.Cells.Item(65, 11).Value = 20678.25
$ws.Cells.Item(65, 13).Value = -17558.25
$ws.Cells.Item(69, 8).Value = 16558.32
$ws.Cells.Item(69, 10).Value = 17299.4
$ws.Cells.Item(69, 12).Value = 51898.2
$ws.Cells.Item(69, 14).Value = -53646.2
$ws.Cells.Item(72, 8).Value = 16558.32
$ws.Cells.Item(72, 10).Value = 17299.4
$ws.Cells.Item(72, 12).Value = 155694.6
$ws.Cells.Item(72, 14).Value = -164430.6
$ws.Cells.Item(86, 8).Value = 12929.125
$ws.Cells.Item(86, 9).Value = 12006.182
$ws.Cells.Item(86, 11).Value = 12006.182
$ws.Cells.Item(86, 13).Value = -10883.182
$ws.Cells.Item(89, 8).Value = 12929.125
$ws.Cells.Item(89, 9).Value = 12006.182
$ws.Cells.Item(89, 11).Value = 60030.91
$ws.Cells.Item(89, 13).Value = -54414.91
$ws.Cells.Item(107, 8).Value = 233.36
$ws.Cells.Item(107, 9).Value = 272.90475
$ws.Cells.Item(107, 11).Value = 272.90475
$ws.Cells.Item(107, 13).Value = 1647.09525
$ws.Cells.Item(111, 8).Value = 252250
$ws.Cells.Item(111, 9).Value = 3000
$ws.Cells.Item(111, 10).Value = 501500
$ws.Cells.Item(111, 11).Value = 9000
$ws.Cells.Item(111, 12).Value = 1504500
$ws.Cells.Item(111, 13).Value = -5933
$ws.Cells.Item(111, 14).Value = -1510634
$ws.Cells.Item(138, 8).Value = 4836.8477
$ws.Cells.Item(138, 9).Value = 5039.4707
$ws.Cells.Item(138, 10).Value = 4262.75
$ws.Cells.Item(138, 11).Value = 15118.4121
$ws.Cells.Item(138, 12).Value = 12788.25
$ws.Cells.Item(138, 13).Value = -9978.4121
$ws.Cells.Item(138, 14).Value = -23068.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7271.7935
$ws.Cells.Item(32, 9).Value = 4029.8718
$ws.Cells.Item(32, 11).Value = 4029.8718
$ws.Cells.Item(32, 13).Value = -3742.8718
$ws.Cells.Item(74, 8).Value = 4624.32
$ws.Cells.Item(74, 9).Value = 2992.3333
$ws.Cells.Item(74, 11).Value = 2992.3333
$ws.Cells.Item(74, 13).Value = -2118.3333
$ws.Cells.Item(77, 8).Value = 4624.32
$ws.Cells.Item(77, 9).Value = 2992.3333
$ws.Cells.Item(77, 11).Value = 14961.6665
$ws.Cells.Item(77, 13).Value = -10593.6665
$ws.Cells.Item(97, 8).Value = 579.29034
$ws.Cells.Item(97, 9).Value = 577.96
$ws.Cells.Item(97, 10).Value = 584.8333
$ws.Cells.Item(97, 11).Value = 577.96
$ws.Cells.Item(97, 12).Value = 584.8333
$ws.Cells.Item(97, 13).Value = -81.96000000000004
$ws.Cells.Item(97, 14).Value = -1576.8333
$ws.Cells.Item(122, 8).Value = 3773.318
$ws.Cells.Item(122, 9).Value = 2525.75
$ws.Cells.Item(122, 10).Value = 16249
$ws.Cells.Item(122, 11).Value = 7577.25
$ws.Cells.Item(122, 12).Value = 48747
$ws.Cells.Item(122, 13).Value = -5127.25
$ws.Cells.Item(122, 14).Value = -53647
$ws.Cells.Item(132, 8).Value = 4041.8572
$ws.Cells.Item(132, 9).Value = 2938.8
$ws.Cells.Item(132, 10).Value = 6799.5
$ws.Cells.Item(132, 11).Value = 8816.400000000001
$ws.Cells.Item(132, 12).Value = 20398.5
$ws.Cells.Item(132, 13).Value = -6286.400000000001
$ws.Cells.Item(132, 14).Value = -25458.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(5, 8).Value = 13575.5
$ws.Cells.Item(5, 9).Value = 1526
$ws.Cells.Item(5, 10).Value = 25625
$ws.Cells.Item(5, 11).Value = 1526
$ws.Cells.Item(5, 12).Value = 25625
$ws.Cells.Item(5, 13).Value = -1413
$ws.Cells.Item(5, 14).Value = -25851
$ws.Cells.Item(20, 8).Value = 3167.5454
$ws.Cells.Item(20, 9).Value = 1872
$ws.Cells.Item(20, 11).Value = 1872
$ws.Cells.Item(20, 13).Value = -1625
$ws.Cells.Item(22, 8).Value = 701.4167
$ws.Cells.Item(22, 9).Value = 510.94116
$ws.Cells.Item(22, 10).Value = 1164
$ws.Cells.Item(22, 11).Value = 510.94116
$ws.Cells.Item(22, 12).Value = 1164
$ws.Cells.Item(22, 13).Value = -337.94116
$ws.Cells.Item(22, 14).Value = -1510
$ws.Cells.Item(107, 8).Value = 1492.1538
$ws.Cells.Item(107, 9).Value = 1514.85
$ws.Cells.Item(107, 10).Value = 1416.5
$ws.Cells.Item(107, 11).Value = 1514.85
$ws.Cells.Item(107, 12).Value = 1416.5
$ws.Cells.Item(107, 13).Value = 405.1500000000001
$ws.Cells.Item(107, 14).Value = -5256.5
$ws.Cells.Item(134, 8).Value = 3151.9565
$ws.Cells.Item(134, 10).Value = 7270.2856
$ws.Cells.Item(134, 12).Value = 21810.8568
$ws.Cells.Item(134, 14).Value = -26880.8568

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 2462.875
$ws.Cells.Item(22, 9).Value = 426.25
$ws.Cells.Item(22, 11).Value = 426.25
$ws.Cells.Item(22, 13).Value = -76.25
$ws.Cells.Item(86, 8).Value = 6274.3335
$ws.Cells.Item(86, 9).Value = 6624.5
$ws.Cells.Item(86, 10).Value = 6099.25
$ws.Cells.Item(86, 11).Value = 6624.5
$ws.Cells.Item(86, 12).Value = 6099.25
$ws.Cells.Item(86, 13).Value = -5501.5
$ws.Cells.Item(86, 14).Value = -8345.25
$ws.Cells.Item(89, 8).Value = 6274.3335
$ws.Cells.Item(89, 9).Value = 6624.5
$ws.Cells.Item(89, 10).Value = 6099.25
$ws.Cells.Item(89, 11).Value = 33122.5
$ws.Cells.Item(89, 12).Value = 30496.25
$ws.Cells.Item(89, 13).Value = -27506.5
$ws.Cells.Item(89, 14).Value = -41728.25
$ws.Cells.Item(105, 8).Value = 2211.1765
$ws.Cells.Item(105, 9).Value = 2224.4375
$ws.Cells.Item(105, 11).Value = 2224.4375
$ws.Cells.Item(105, 13).Value = -477.4375
$ws.Cells.Item(122, 8).Value = 2923.6155
$ws.Cells.Item(122, 9).Value = 2430.2
$ws.Cells.Item(122, 10).Value = 3596.4546
$ws.Cells.Item(122, 11).Value = 7290.599999999999
$ws.Cells.Item(122, 12).Value = 10789.3638
$ws.Cells.Item(122, 13).Value = -4840.599999999999
$ws.Cells.Item(122, 14).Value = -15689.3638
$ws.Cells.Item(134, 8).Value = 4052.5264
$ws.Cells.Item(134, 9).Value = 3116.5715
$ws.Cells.Item(134, 10).Value = 6673.2
$ws.Cells.Item(134, 11).Value = 9349.7145
$ws.Cells.Item(134, 12).Value = 20019.6
$ws.Cells.Item(134, 13).Value = -6814.7145
$ws.Cells.Item(134, 14).Value = -25089.6

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(15, 8).Value = 61.363636
$ws.Cells.Item(15, 9).Value = 48.125
$ws.Cells.Item(15, 10).Value = 96.666664
$ws.Cells.Item(15, 11).Value = 144.375
$ws.Cells.Item(15, 12).Value = 289.999992
$ws.Cells.Item(15, 13).Value = -4.375
$ws.Cells.Item(15, 14).Value = -569.999992
$ws.Cells.Item(55, 8).Value = 6299.5
$ws.Cells.Item(55, 10).Value = 8232.666999999999
$ws.Cells.Item(55, 12).Value = 24698.001
$ws.Cells.Item(55, 14).Value = -25052.001
$ws.Cells.Item(122, 8).Value = 0
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 2048.3333
$ws.Cells.Item(132, 9).Value = 1992.75
$ws.Cells.Item(132, 11).Value = 17934.75
$ws.Cells.Item(132, 13).Value = -15404.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 6218.0293
$ws.Cells.Item(132, 9).Value = 5368.5386
$ws.Cells.Item(132, 11).Value = 16105.6158
$ws.Cells.Item(132, 13).Value = -13575.6158
$ws.Cells.Item(136, 8).Value = 20923.5
$ws.Cells.Item(136, 10).Value = 20923.5
$ws.Cells.Item(136, 12).Value = 62770.5
$ws.Cells.Item(136, 14).Value = -67870.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 6182.9585
$ws.Cells.Item(61, 9).Value = 2699.3076
$ws.Cells.Item(61, 10).Value = 10300
$ws.Cells.Item(61, 11).Value = 2699.3076
$ws.Cells.Item(61, 12).Value = 10300
$ws.Cells.Item(61, 13).Value = -2497.3076
$ws.Cells.Item(61, 14).Value = -10704
$ws.Cells.Item(100, 8).Value = 4445.1
$ws.Cells.Item(100, 9).Value = 3827.4443
$ws.Cells.Item(100, 11).Value = 3827.4443
$ws.Cells.Item(100, 13).Value = -3286.4443
$ws.Cells.Item(113, 8).Value = 6182.9585
$ws.Cells.Item(113, 9).Value = 2699.3076
$ws.Cells.Item(113, 10).Value = 10300
$ws.Cells.Item(113, 11).Value = 2699.3076
$ws.Cells.Item(113, 12).Value = 10300
$ws.Cells.Item(113, 13).Value = -529.3076000000001
$ws.Cells.Item(113, 14).Value = -14640
$ws.Cells.Item(136, 8).Value = 7447.113
$ws.Cells.Item(136, 9).Value = 4974.2666
$ws.Cells.Item(136, 11).Value = 14922.7998
$ws.Cells.Item(136, 13).Value = -12372.7998

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 14).ClearContents()
$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 14).ClearContents()
$ws.Cells.Item(107, 8).Value = 2764.6667
$ws.Cells.Item(107, 9).Value = 2577.7273
$ws.Cells.Item(107, 11).Value = 7733.1819
$ws.Cells.Item(107, 13).Value = -5813.1819
$ws.Cells.Item(136, 8).Value = 3732.85
$ws.Cells.Item(136, 9).Value = 3013.1455
$ws.Cells.Item(136, 10).Value = 11649.6
$ws.Cells.Item(136, 11).Value = 9039.4365
$ws.Cells.Item(136, 12).Value = 34948.8
$ws.Cells.Item(136, 13).Value = -6489.4365
$ws.Cells.Item(136, 14).Value = -40048.8
